# Generate Report for Handback
#
# The second tracked source file (e6052173-41f7-4b7d-8df8-c8821d2e3359.md)
# has completed localization round-trip for both target locales, so its
# status flips from "Ready for handoff" to "Handed back: in sync with en-US"
# and its "Latest Handback DateTime" is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: both locale status columns for that file's row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: status + handback datetime for that file's row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("G3").Value = "2016-02-24 07:10:33"

# --- de-de sheet: status + handback datetime for that file's row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusHandedBack
$dede.Range("G3").Value = "2016-02-24 07:10:57"
